$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,44.29505033333334,132.885151,0.9830698162761968,0.9830698162761969,3,1,58.95713633333333,176.871409,0.4863146960083892,0.4863146960083893,2611.509321394195,23503.58389254776,0.4780812988573817,0.4780812988573818),
    @(3,1,44.29505033333334,132.885151,0.9830698162761968,0.9830698162761969,3,1,8.99153,26.97459,0.07416766570679004,0.07416766570679005,398.2802739236766,3584.52246531309,0.07291199350000846,0.07291199350000849),
    @(3,1,44.29505033333334,132.885151,0.9830698162761968,0.9830698162761969,3,1,42.51661933333333,127.549858,0.3507032073181665,0.3507032073181665,1883.275793373173,16949.48214035856,0.3447657375857429,0.344765737585743),
    @(3,1,44.29505033333334,132.885151,0.9830698162761968,0.9830698162761969,3,1,10.76719366666667,32.301581,0.0888144309666542,0.08881443096665421,476.9333854137479,4292.400468723731,0.08731078633306372,0.08731078633306373),
    @(1,0.3333333333333333,0.307031,0.921093,0.006814145293655052,0.006814145293655053,3,1,58.95713633333333,176.871409,0.4863146960083892,0.4863146960083893,18.10166852555967,162.915016730037,0.003313818997040853,0.003313818997040854),
    @(1,0.3333333333333333,0.307031,0.921093,0.006814145293655052,0.006814145293655053,3,1,8.99153,26.97459,0.07416766570679004,0.07416766570679005,2.76067844743,24.84610602687,0.0005053892502173045,0.0005053892502173047),
    @(1,0.3333333333333333,0.307031,0.921093,0.006814145293655052,0.006814145293655053,3,1,42.51661933333333,127.549858,0.3507032073181665,0.3507032073181665,13.05392015053267,117.485281354794,0.002389742609616816,0.002389742609616817),
    @(1,0.3333333333333333,0.307031,0.921093,0.006814145293655052,0.006814145293655053,3,1,10.76719366666667,32.301581,0.0888144309666542,0.08881443096665421,3.305862238670333,29.752760148033,0.0006051944367800782,0.0006051944367800784),
    @(3,1,0.3685326666666667,1.105598,0.008179093108268589,0.008179093108268589,3,1,58.95713633333333,176.871409,0.4863146960083892,0.4863146960083893,21.72763067195356,195.548676047582,0.00397761317857195,0.00397761317857195),
    @(3,1,0.3685326666666667,1.105598,0.008179093108268589,0.008179093108268589,3,1,8.99153,26.97459,0.07416766570679004,0.07416766570679005,3.313672528313333,29.82305275482,0.0006066242434387749,0.000606624243438775),
    @(3,1,0.3685326666666667,1.105598,0.008179093108268589,0.008179093108268589,3,1,42.51661933333333,127.549858,0.3507032073181665,0.3507032073181665,15.66876310056489,141.018867905084,0.002868434186023705,0.002868434186023706),
    @(3,1,0.3685326666666667,1.105598,0.008179093108268589,0.008179093108268589,3,1,10.76719366666667,32.301581,0.0888144309666542,0.08881443096665421,3.968062594493111,35.712563350438,0.0007264215002341577,0.0007264215002341578),
    @(2,0.6666666666666666,0.08727466666666667,0.261824,0.001936945321879485,0.001936945321879485,3,1,58.95713633333333,176.871409,0.4863146960083892,0.4863146960083893,5.145464421112889,46.309179790016,0.0009419649753946933,0.0009419649753946935),
    @(2,0.6666666666666666,0.08727466666666667,0.261824,0.001936945321879485,0.001936945321879485,3,1,8.99153,26.97459,0.07416766570679004,0.07416766570679005,0.7847327835733333,7.06259505216,0.0001436587131254885,0.0001436587131254885),
    @(2,0.6666666666666666,0.08727466666666667,0.261824,0.001936945321879485,0.001936945321879485,3,1,42.51661933333333,127.549858,0.3507032073181665,0.3507032073181665,3.710623780110222,33.395614020992,0.0006792929367830536,0.0006792929367830539),
    @(2,0.6666666666666666,0.08727466666666667,0.261824,0.001936945321879485,0.001936945321879485,3,1,10.76719366666667,32.301581,0.0888144309666542,0.08881443096665421,0.9397032381937777,8.457329143744,0.0001720286965762493,0.0001720286965762494)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 5).Value = $row[$j]
    }
}

Write-Output "done"